$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update only the cells whose values changed, per the target diff.
$ws.Cells.Item(2, 2).Value = '航天发展'
$ws.Cells.Item(2, 3).Value = '航天发展'

$ws.Cells.Item(3, 1).Value = '航天动力'
$ws.Cells.Item(3, 2).Value = 'N摩尔-U'
$ws.Cells.Item(3, 3).Value = '海王生物'

$ws.Cells.Item(4, 2).Value = '龙洲股份'

$ws.Cells.Item(5, 1).Value = '顺灏股份'
$ws.Cells.Item(5, 2).Value = '海王生物'
$ws.Cells.Item(5, 3).Value = '摩尔线程'

$ws.Cells.Item(6, 1).Value = '龙洲股份'
$ws.Cells.Item(6, 2).Value = '实达集团'
$ws.Cells.Item(6, 3).Value = '实达集团'

$ws.Cells.Item(7, 1).Value = '平潭发展'
$ws.Cells.Item(7, 2).Value = '顺灏股份'
$ws.Cells.Item(7, 3).Value = '龙洲股份'

$ws.Cells.Item(8, 1).Value = 'N摩尔-U'
$ws.Cells.Item(8, 2).Value = '赢时胜'
$ws.Cells.Item(8, 3).Value = '合富中国'

$ws.Cells.Item(9, 2).Value = '航天机电'
$ws.Cells.Item(9, 3).Value = '航天动力'

$ws.Cells.Item(10, 1).Value = '实达集团'
$ws.Cells.Item(10, 2).Value = '合富中国'
$ws.Cells.Item(10, 3).Value = '顺灏股份'

$ws.Cells.Item(11, 1).Value = '永鼎股份'
$ws.Cells.Item(11, 2).Value = '航天动力'
$ws.Cells.Item(11, 3).Value = '海欣食品'

$ws.Cells.Item(12, 2).Value = '航天科技'
$ws.Cells.Item(12, 3).Value = '太阳电缆'

$ws.Cells.Item(13, 1).Value = '航天机电'
$ws.Cells.Item(13, 2).Value = '东方财富'
$ws.Cells.Item(13, 3).Value = '国机重装'

$ws.Cells.Item(14, 1).Value = '海欣食品'
$ws.Cells.Item(14, 2).Value = '平潭发展'
$ws.Cells.Item(14, 3).Value = '航天机电'

$ws.Cells.Item(15, 1).Value = '航天科技'
$ws.Cells.Item(15, 3).Value = '瑞康医药'

$ws.Cells.Item(16, 1).Value = '国机重装'
$ws.Cells.Item(16, 2).Value = '国机重装'
$ws.Cells.Item(16, 3).Value = '通宇通讯'

$ws.Cells.Item(17, 1).Value = '太阳电缆'
$ws.Cells.Item(17, 2).Value = '永鼎股份'
$ws.Cells.Item(17, 3).Value = '特发信息'

$ws.Cells.Item(18, 1).Value = '瑞康医药'
$ws.Cells.Item(18, 2).Value = '瑞康医药'
$ws.Cells.Item(18, 3).Value = '雷科防务'

$ws.Cells.Item(19, 1).Value = '赢时胜'
$ws.Cells.Item(19, 2).Value = '宁波韵升'
$ws.Cells.Item(19, 3).Value = '安记食品'

$ws.Cells.Item(20, 1).Value = '特发信息'
$ws.Cells.Item(20, 2).Value = '太阳电缆'
$ws.Cells.Item(20, 3).Value = '雪人集团'

$ws.Cells.Item(21, 1).Value = '上海瀚讯'
$ws.Cells.Item(21, 2).Value = '特发信息'
$ws.Cells.Item(21, 3).Value = '双星新材'
